$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# Cell B27 corresponds to key "columnsPartnershipU1a" (row 27: A27="columnsPartnershipU1a").
# The value is changed from the number 30 to the text value "28" entered with a leading
# apostrophe (quote-prefix) in Excel, so it is stored as text rather than as a number,
# using the font from the "VALUE" style (fontId 1).
$cell = $ws.Range("B27")
$cell.Value = "'28"

# Update the active selection to match the edited cell, as happens naturally in Excel
# after editing B27.
$ws.Range("B28").Select()
